$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.822.84"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").Value = "1.965.24"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.37"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4778"
$ws.Range("E7").Value = "  -4.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4048"
$ws.Range("E8").Value = "  -4.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.74"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08504"
$ws.Range("E10").Value = "  -5.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.062"
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.49"
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("D13").Value = "1.974.64"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.665"
$ws.Range("E14").Value = "  -5.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.254"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.015"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.05"
$ws.Range("E17").Value = "  -4.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001070"
$ws.Range("E18").Value = "  -3.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06617"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.71"
$ws.Range("E20").Value = "  -5.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.791"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").Value = "28.837.10"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.58"
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.290"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "2.182.74"
$ws.Range("E26").Value = "  -3.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.54"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.25"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.970"
$ws.Range("E29").Value = "  -6.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.165"
$ws.Range("E30").Value = "  -6.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.34"
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.007"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09608"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.467"
$ws.Range("E34").Value = "  -6.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.681"
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.680"
$ws.Range("E36").Value = "  -3.33%  "
$ws.Range("E37").Value = "  -4.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.271"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.801"
$ws.Range("E39").Value = "  -6.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06215"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6253"
$ws.Range("E41").Value = "  -4.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.13"
$ws.Range("E42").Value = "  -4.83%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1920"
$ws.Range("E44").Value = "  -6.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.345"
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5978"
$ws.Range("E46").Value = "  -5.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.06"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.085"
$ws.Range("E48").Value = "  -5.55%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000337"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("B50").Value = "PancakeSwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.420"
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06843"
$ws.Range("E51").Value = "  -2.03%  "
